# Update existing rows 2-29 to reflect the refreshed weekly price feed
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "2023-09-22"
$ws.Range("M2").Value = 30
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 21000
$ws.Range("S2").Value = 2100
$ws.Range("D3").Value = "2023-09-11"
$ws.Range("M3").Value = 40
$ws.Range("N3").Value = 22000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 22000
$ws.Range("S3").Value = 2200
$ws.Range("D4").Value = "2023-09-20"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 22000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 22000
$ws.Range("S4").Value = 2200
$ws.Range("D5").Value = "2023-09-25"
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 22000
$ws.Range("O5").Value = 22000
$ws.Range("P5").Value = 22000
$ws.Range("S5").Value = 2200
$ws.Range("D6").Value = "2023-09-07"
$ws.Range("N6").Value = 22000
$ws.Range("O6").Value = 22000
$ws.Range("P6").Value = 22000
$ws.Range("S6").Value = 2200
$ws.Range("D7").Value = "2021-09-07"
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = 21000
$ws.Range("P7").Value = 21500
$ws.Range("S7").Value = 2150
$ws.Range("D8").Value = "2021-09-22"
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 31000
$ws.Range("O8").Value = 32000
$ws.Range("P8").Value = 31500
$ws.Range("S8").Value = 3150
$ws.Range("D9").Value = "2021-09-22"
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 30000
$ws.Range("O9").Value = 30000
$ws.Range("P9").Value = 30000
$ws.Range("S9").Value = 3000
$ws.Range("D10").Value = "2022-10-07"
$ws.Range("M10").Value = 60
$ws.Range("O10").Value = 24000
$ws.Range("P10").Value = 23500
$ws.Range("S10").Value = 2350
$ws.Range("D11").Value = "2023-10-13"
$ws.Range("M11").Value = 80
$ws.Range("N11").Value = 22000
$ws.Range("P11").Value = 22000
$ws.Range("S11").Value = 2200
$ws.Range("D12").Value = "2023-10-10"
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = 22000
$ws.Range("O12").Value = 22000
$ws.Range("P12").Value = 22000
$ws.Range("S12").Value = 2200
$ws.Range("D13").Value = "2023-09-27"
$ws.Range("N13").Value = 23000
$ws.Range("O13").Value = 23000
$ws.Range("P13").Value = 23000
$ws.Range("S13").Value = 2300
$ws.Range("D14").Value = "2023-10-18"
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 30
$ws.Range("N14").Value = 21000
$ws.Range("O14").Value = 21000
$ws.Range("P14").Value = 21000
$ws.Range("S14").Value = 2100
$ws.Range("D15").Value = "2023-10-06"
$ws.Range("N15").Value = 22000
$ws.Range("O15").Value = 22000
$ws.Range("P15").Value = 22000
$ws.Range("S15").Value = 2200
$ws.Range("D16").Value = "2021-09-13"
$ws.Range("D17").Value = "2022-10-05"
$ws.Range("M17").Value = 120
$ws.Range("N17").Value = 25000
$ws.Range("O17").Value = 26000
$ws.Range("P17").Value = 25500
$ws.Range("S17").Value = 2550
$ws.Range("D18").Value = "2022-10-14"
$ws.Range("L18").Value = "Especial"
$ws.Range("M18").Value = 60
$ws.Range("N18").Value = 24000
$ws.Range("O18").Value = 25000
$ws.Range("P18").Value = 24500
$ws.Range("S18").Value = 2450
$ws.Range("D19").Value = "2022-10-14"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 120
$ws.Range("N19").Value = 21000
$ws.Range("O19").Value = 22000
$ws.Range("P19").Value = 21500
$ws.Range("S19").Value = 2150
$ws.Range("D20").Value = "2021-09-09"
$ws.Range("M20").Value = 60
$ws.Range("N20").Value = 21000
$ws.Range("O20").Value = 22000
$ws.Range("P20").Value = 21500
$ws.Range("S20").Value = 2150
$ws.Range("D21").Value = "2021-09-08"
$ws.Range("M21").Value = 60
$ws.Range("N21").Value = 21000
$ws.Range("P21").Value = 21500
$ws.Range("S21").Value = 2150
$ws.Range("D22").Value = "2023-10-04"
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 30
$ws.Range("N22").Value = 21000
$ws.Range("O22").Value = 21000
$ws.Range("P22").Value = 21000
$ws.Range("S22").Value = 2100
$ws.Range("D23").Value = "2021-10-18"
$ws.Range("M23").Value = 30
$ws.Range("N23").Value = 23000
$ws.Range("O23").Value = 24000
$ws.Range("P23").Value = 23500
$ws.Range("S23").Value = 2350
$ws.Range("D24").Value = "2022-11-09"
$ws.Range("L24").Value = "Especial"
$ws.Range("M24").Value = 30
$ws.Range("N24").Value = 25000
$ws.Range("O24").Value = 25000
$ws.Range("P24").Value = 25000
$ws.Range("S24").Value = 2500
$ws.Range("D25").Value = "2022-11-09"
$ws.Range("M25").Value = 80
$ws.Range("N25").Value = 23000
$ws.Range("O25").Value = 24000
$ws.Range("P25").Value = 23500
$ws.Range("S25").Value = 2350
$ws.Range("D26").Value = "2023-09-04"
$ws.Range("D27").Value = "2022-11-03"
$ws.Range("L27").Value = "Especial"
$ws.Range("M27").Value = 60
$ws.Range("N27").Value = 26000
$ws.Range("O27").Value = 26000
$ws.Range("P27").Value = 26000
$ws.Range("S27").Value = 2600
$ws.Range("D28").Value = "2023-10-17"
$ws.Range("M28").Value = 60
$ws.Range("N28").Value = 21000
$ws.Range("O28").Value = 21000
$ws.Range("P28").Value = 21000
$ws.Range("S28").Value = 2100
$ws.Range("D29").Value = "2021-09-21"
$ws.Range("L29").Value = "Especial"
$ws.Range("M29").Value = 60
$ws.Range("N29").Value = 31000
$ws.Range("O29").Value = 32000
$ws.Range("P29").Value = 31500
$ws.Range("S29").Value = 3150

# Append the new weekly record as row 30
$ws.Range("D30").NumberFormat = $ws.Range("D29").NumberFormat
$ws.Range("A30").Value = 7
$ws.Range("B30").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C30").Value = "Ñuble"
$ws.Range("D30").Value = "2021-09-21"
$ws.Range("E30").Value = 16
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100107
$ws.Range("H30").Value = "Otros"
$ws.Range("I30").Value = 100107002
$ws.Range("J30").Value = "Chirimoya"
$ws.Range("K30").Value = "Cultivar IV Región"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 30
$ws.Range("N30").Value = 30000
$ws.Range("O30").Value = 30000
$ws.Range("P30").Value = 30000
$ws.Range("Q30").Value = "`$/bandeja 10 kilos"
$ws.Range("R30").Value = "Provincia de Limarí"
$ws.Range("S30").Value = 3000
$ws.Range("T30").Value = 10

Write-Host "Updated rows 2-29 and appended new row 30"
